# The commit swaps the two embedded theme parts: the "Integral" theme
# (ppt/theme/theme1.xml, used by the real slide master) and the
# "Office Theme" (ppt/theme/theme2.xml, used only by the notes master)
# change places - the slide master now uses the stock "Office Theme"
# color palette instead of "Integral".
#
# The font scheme and format scheme (fills/lines/effects) are identical
# between the two themes already, so the only thing that visibly moves
# is the 12-slot color scheme. Apply the Office Theme palette to the
# presentation's color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink, in that order) via the slide's ColorScheme, which is backed
# by the slide master's theme part.

$p = $ppt.ActivePresentation

$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$s = $p.Slides.Item(1)
$cs = $s.ColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $rgbLong = $r + ($g * 256) + ($b * 65536)
    $cs.Colors($i).RGB = $rgbLong
}
